# Edit script: applies three paragraph-level text revisions to the coop work
# journal, matching the target OOXML diff:
#  1. Split "in" out into its own run (wrapped in w:proofErr gramStart/gramEnd)
#     inside the 2024-05-08 paragraph mentioning "another of Sean's training
#     sessions".
#  2. Split "thorough" out into its own run (wrapped in w:proofErr
#     gramStart/gramEnd) inside the 2024-05-09 "advanced report training"
#     paragraph.
#  3. Rewrite the 2024-05-15 entry (collapsing its many runs into two) and
#     append a new 2024-05-16 journal entry (blank paragraph + bold date
#     heading + body paragraph).
#
# Because the edits splice in <w:proofErr/> elements and brand-new
# paragraphs that plain Find/Replace can't express, each change is applied
# by locating the target paragraph with Range.Find and then replacing that
# paragraph's full contents with literal OOXML via Range.InsertXML (which
# replaces the exact range it is called on).

$d = $word.ActiveDocument

function Get-ParagraphRangeByText($doc, [string]$needle) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $needle"
    }
    return $rng.Paragraphs(1).Range
}

# --- Change 1: 2024-05-08 entry, "in another of Sean's training sessions" ---
$para1 = Get-ParagraphRangeByText $d "Tuned in to another one of Sean's training sessions"
$para1.InsertXML('<w:p w14:paraId="08BDAB22" w14:textId="0C3F777F" w:rsidR="00EC205B" w:rsidRDefault="00EC205B" w:rsidP="00EA492B"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00EC205B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Tuned in to another one of Sean''s training sessions. This one was a bit of a slog, 4 hours long total. However, this session was with a new group of people and this group is much more engaged than the previous. They were asking questions following along well which made the 4 hours pass a bit quicker. I mostly continued with reports during this time, so the meeting served as background noise. Initially the company who employed us wanted this new data system to go live </w:t></w:r><w:r w:rsidR="00810420" w:rsidRPr="00EC205B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>from</w:t></w:r><w:r w:rsidRPr="00EC205B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the end of May to the beginning of June. However, the lead on this project at our employer said that the yearly shutdown they do in September may be </w:t></w:r><w:r w:rsidRPr="00EC205B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">moved to June instead. This would give us a little more time to flesh out the system and polish some aspects. We are still operating as if the timeline is the same but will take advantage of the extra time if the decision to change the shutdown ends up being made. Afterwards, I was also </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>in</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> another of Sean''s training sessions, but this one was only an hour. It was mostly showing off the trend viewer client app, query tool and excel template that I made yesterday. This session went well too, and Sean got a lot of positive feedback which he appreciated.</w:t></w:r></w:p>')

# --- Change 2: 2024-05-09 "advanced report training" entry, "thorough" ---
$para2 = Get-ParagraphRangeByText $d "This morning leads off with an advanced report training"
$para2.InsertXML('<w:p w14:paraId="68DE983D" w14:textId="4BD7A878" w:rsidR="000E6ED2" w:rsidRDefault="00584DDC" w:rsidP="00EA492B"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00584DDC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This morning leads off with an advanced report training. This initial group was less interesting than yesterday’s. They did not interact as much as the last group. This session basically just covered charts in the report studio. The charts are like Excel charts but with a little more convoluted feature. This training went smoothly for the most part, sometimes the report studio has errors when you try and make object near the edge of the report page causes you to use some less than convenient work arounds. New to this training session, Sean let me teach the report template section. I was nervous at first but eventually just focused on doing a good job and forgot about being nervous. Post training meeting both the lead at the company employing us and Sean complimented how clearly and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>thorough</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I was explaining the report template features, which made me feel a lot more confident in my ability. Then for the rest of the day I worked on fleshing out the engineering unit standardization I was working on yesterday and finished up the document and gave it to Sean. Researching symbols and standards for units made me realize how inconsistent they </w:t></w:r><w:r w:rsidRPr="00584DDC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>can be sometimes despite having a standard. Regardless I finish out the last hour of the day by doing some more reports and heading home for the weekend.</w:t></w:r></w:p>')

# --- Change 3: 2024-05-15 entry rewrite + new 2024-05-16 entry appended ---
$para3 = Get-ParagraphRangeByText $d "Today we took another big step forward to the go live"
$para3.InsertXML('<w:p w14:paraId="7E5C231C" w14:textId="19164721" w:rsidR="00810420" w:rsidRPr="00EA492B" w:rsidRDefault="00810420" w:rsidP="00810420"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00810420"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Today we took another big step forward to the go live for the new system for our employer. Since most of the reports are finished and compiled, we now must vet their functionality. This morning Sean asked me to make a sheet to track and test all the reports document individual issues with each of their tables or views. This is pretty much going through all the </w:t></w:r><w:r w:rsidRPr="00810420"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>reports again, like I have been since starting. I''m quite proud of the sheet I made. It works decently well and modular which will allow for updates. It seems like a lot of the tags have broken values or the way they are being queried is incorrect. I say this because a lot of the tables have values displaying "N/A" or "Err" meaning no data or there are null values respectively. This is a problem likely with the tag itself, somewhere along the retrieval line, there could be a typo with the tag or something along those lines. For now, we''re just taking note of them and then going to return and amend them later. Hopefully, when we go back to fix them, the fixes will be as simple as fixing a type (although it''s **never** usually that easy). So, another long but relatively simple task that will likely close out the week for me.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>2024-05-16</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Mostly report verification today, going through the standard tests and verifying everything is in working order. Made it through a good </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>junk</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> report </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>today</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. However, it is a little annoying because Sean is doing training sessions at the same time I''m using the web view for the reports. Since he is doing training, every so often he </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> reload the whole project containing all the reports to update changes he makes to his examples. This inadvertently </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>caused</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> me to be logged out of the online portal due to the web view being updated. This causes small interruptions when I''m trying to be efficient, having to log back in causing some productivity halting. But today, a bigger interruption was experienced. Towards the afternoon, the whole project housing the reports crashed and I was no longer able to get onto the report studio software. I got in contact with the IT department for our employer and he basically said he had to go physically </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>to find</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> out what was going on. Later, Sean messaged me saying that the whole server needed to be restarted before I </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>could</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> continue. So pretty much no work got done for the last two hours of today.</w:t></w:r></w:p>')

Write-Output "done"
